$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-09 Friday" "2024-08-10 Saturday"

Replace-Text "830÷9=" "852÷2="
Replace-Text "367÷7=" "525÷5="
Replace-Text "517÷5=" "890÷6="
Replace-Text "293÷5=" "486÷5="
Replace-Text "220÷8=" "647÷8="
Replace-Text "806÷8=" "898÷9="
Replace-Text "692÷5=" "245÷5="
Replace-Text "646÷9=" "330÷4="
Replace-Text "258÷2=" "706÷5="
Replace-Text "321÷6=" "586÷5="
Replace-Text "269÷9=" "265÷9="
Replace-Text "244÷9=" "439÷4="
Replace-Text "748÷7=" "346÷8="
Replace-Text "785÷3=" "605÷6="
Replace-Text "179÷2=" "424÷2="
Replace-Text "505÷9=" "261÷7="
Replace-Text "836÷9=" "829÷5="
Replace-Text "129÷2=" "200÷6="
Replace-Text "735÷6=" "948÷7="
Replace-Text "459÷8=" "809÷2="
Replace-Text "646÷8=" "206÷4="
Replace-Text "517÷3=" "812÷8="
Replace-Text "295÷9=" "800÷4="
Replace-Text "934÷5=" "433÷6="
Replace-Text "186÷4=" "826÷5="
